$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 29   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# --- Simple numeric value updates (rows 15-30 crime stats) ---
$values = @{
    "F15" = 1
    "H15" = -50
    "N15" = 100
    "C16" = 1
    "D16" = 1
    "E16" = 0
    "F16" = 11
    "G16" = 8
    "H16" = 37.5
    "I16" = 165
    "J16" = 99
    "K16" = 66.666666666666
    "L16" = 60.194174757281
    "M16" = 1.851851851851
    "N16" = -80.858468677494
    "C17" = 6
    "E17" = 100
    "F17" = 15
    "G17" = 16
    "H17" = -6.25
    "I17" = 191
    "J17" = 178
    "K17" = 7.303370786516
    "L17" = 45.801526717557
    "M17" = 78.504672897196
    "N17" = -28.464419475655
    "C18" = 1
    "D18" = 4
    "E18" = -75
    "F18" = 14
    "H18" = -17.647058823529
    "I18" = 141
    "J18" = 155
    "K18" = -9.032258064516
    "L18" = -4.081632653061
    "M18" = -32.211538461538
    "N18" = -89.228418640183
    "C19" = 10
    "D19" = 11
    "E19" = -9.090909090909
    "F19" = 51
    "H19" = 8.510638297872
    "I19" = 597
    "J19" = 381
    "K19" = 56.692913385826
    "L19" = 68.644067796610
    "M19" = 51.139240506329
    "N19" = -23.949044585987
    "C20" = 6
    "D20" = 2
    "E20" = 200
    "G20" = 16
    "H20" = 37.5
    "I20" = 202
    "J20" = 157
    "K20" = 28.662420382165
    "L20" = 36.486486486486
    "M20" = 13.483146067415
    "N20" = -88.683473389355
    "C21" = 24
    "D21" = 21
    "E21" = 14.285714285714
    "F21" = 114
    "G21" = 106
    "H21" = 7.547169811320
    "I21" = 1315
    "J21" = 983
    "K21" = 33.774160732451
    "L21" = 47.091722595078
    "M21" = 23.474178403755
    "N21" = -73.830845771144
    "C22" = 3
    "F22" = 10
    "G22" = 2
    "H22" = 400
    "I22" = 71
    "K22" = 108.823529411765
    "L22" = 136.666666666667
    "M22" = 61.363636363636
    "C24" = 27
    "D24" = 21
    "E24" = 28.571428571428
    "F24" = 131
    "G24" = 103
    "H24" = 27.184466019417
    "I24" = 1282
    "J24" = 1086
    "K24" = 18.047882136279
    "L24" = 65.206185567010
    "M24" = 58.076448828606
    "C25" = 13
    "D25" = 9
    "E25" = 44.444444444444
    "F25" = 59
    "G25" = 37
    "H25" = 59.459459459459
    "I25" = 485
    "J25" = 405
    "K25" = 19.753086419753
    "L25" = 68.989547038327
    "M25" = 12.268518518518
    "F26" = 2
    "H26" = 0
    "C27" = 2
    "D27" = 1
    "E27" = 100
    "I27" = 78
    "J27" = 53
    "K27" = 47.169811320754
    "L27" = 39.285714285714
    "J30" = 12
    "K30" = -25
}
foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}

# --- Special cells: numeric <-> text("0"/"***.*")  type changes ---
# Use stable donor cells (row 23, untouched by this edit) to copy both
# the shared-string value and the style atomically.
$ws.Range("C23").Copy($ws.Range("C15"))   # 1 -> "0" (text)
$ws.Range("C23").Copy($ws.Range("D22"))   # 1 -> "0" (text)
$ws.Range("E23").Copy($ws.Range("E22"))   # 100 -> "***.*" (text)
$ws.Range("C23").Copy($ws.Range("C26"))   # 1 -> "0" (text)

# text -> numeric: paste the number-style format from a donor, then set value
$ws.Range("I30").Copy()
$ws.Range("D30").PasteSpecial(-4122)      # xlPasteFormats (style 15)
$ws.Range("G30").PasteSpecial(-4122)      # xlPasteFormats (style 15)
$ws.Range("K30").Copy()
$ws.Range("E30").PasteSpecial(-4122)      # xlPasteFormats (style 16)
$ws.Range("H30").PasteSpecial(-4122)      # xlPasteFormats (style 16)
$excel.CutCopyMode = $false

$ws.Range("D30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("H30").Value = -100